# Add two new columns (I, J) with header labels "I0" / "IF" and three
# data rows of numbers, matching the existing sheet layout (header row 1,
# data rows 2-4). The existing header cells (B1:H1) and the row-index
# column (A2:A4) use a shared bold/bordered/centered style, which new
# header cells I1/J1 should also pick up; data cells I2:J4 stay unstyled,
# same as the other numeric data cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting from an existing header cell (H1) onto the two new
# header cells so they match the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New data cells for rows 2-4.
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 7

$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 7
